$wb = $excel.ActiveWorkbook

# Duplicate the "Spain" sheet (same layout/styles/values) and place the
# copy immediately after it - this becomes the new "Turkey" sheet.
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Turkey only has two loading rows (the "MZX254" panel row used by Spain
# is not applicable), so drop data row 9 - this shifts former row 10
# (P405D) up to row 9 and updates the sheet's used-range automatically.
$turkey.Rows(9).Delete()

# Point the "User Story" reference (B4) at the Turkey-specific ticket.
$turkey.Range("B4").Value = "NGC-3191/T3331/T3332/T3330"

# Spain is no longer the active tab; its selection reverts to the whole
# used range (matching the recorded sheetView for that sheet).
$spain.Range("A1:U10").Select()

# Turkey ends up as the active/selected tab, with its recorded selection.
$turkey.Activate()
$turkey.Range("G10").Select()
